# Automatic hashcode update
# Updates the hashcode values (column B) for a set of rows identified by
# their row number in the "hashcode.csv" worksheet, matching the upstream
# diff applied to data/metadata/hashcode.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "1352d9b99bf06626ff80952eda02d7d2"
$ws.Range("B44").Value = "775da89266fde57dfe7ca7c89abf5d91"
$ws.Range("B74").Value = "8a74666dc4ebb183229cedc771aa374f"
$ws.Range("B89").Value = "e5a9c26e094a5557ae9c4aa83e416d55"
$ws.Range("B99").Value = "0c473cacc596f7b80f753639d0d0ca9c"
$ws.Range("B110").Value = "8c9098805d070995ea6995c660cc73a1"
$ws.Range("B121").Value = "27c1bb70cb640d5ca20a759347c927c8"
$ws.Range("B126").Value = "d10e7f81f334a0777fde493a9d8be49f"
$ws.Range("B154").Value = "e9828e955ed4896624069e2230da5da2"
$ws.Range("B160").Value = "f3de5288eeaf606f566c40f38f1f948a"
$ws.Range("B168").Value = "bc95cae257a5ff8399d8aa38ac0096e0"
$ws.Range("B191").Value = "c73e5ad0a567948972aa3db3a087d497"
$ws.Range("B246").Value = "63e662a7fda656ea7e5a89ae489194ed"
$ws.Range("B276").Value = "aa67a5307aa39ccf124b7bd08af4cef5"
$ws.Range("B278").Value = "9283cf6e227051ed64790cd8214746ac"
$ws.Range("B293").Value = "21201fdc44ce87e98d9209da669acf6b"
$ws.Range("B302").Value = "86f7fce844f6cb8603cc25fce67ebddc"
$ws.Range("B345").Value = "3d3502f758d76be92c0f4e2ea3201dd1"
$ws.Range("B410").Value = "c099fb6691a99c8e052f7a2590abc91e"
$ws.Range("B446").Value = "96182ee2b15c2bf273c450cd40d32591"
$ws.Range("B460").Value = "dcb2ec38d53f4f46da20eaa229beb33b"
$ws.Range("B480").Value = "f23b3dca7b162c63f81a3379142179f4"
$ws.Range("B500").Value = "f359e34c0328c91d3de985593087f5b0"
$ws.Range("B501").Value = "4d6e74117798826934f5aaa2f340fd95"
$ws.Range("B517").Value = "cd14a256e4239dd10d8a16192838843e"
$ws.Range("B534").Value = "76da3783aa2a61aa6867b6ba825b3179"
$ws.Range("B547").Value = "61c4f18193adac7d146bc75c0f680430"
$ws.Range("B550").Value = "764c658498c1acf6a3e233b45ec55287"
$ws.Range("B553").Value = "58d85ba2051dd71507a5e4255d2e5b94"
$ws.Range("B566").Value = "dbea5a0e4f8a16f2177f6d333e483de5"
$ws.Range("B665").Value = "4623493d74b0998e011a5f81554979a9"
$ws.Range("B756").Value = "d433d8485854eb32ee18383271ee08ce"
$ws.Range("B761").Value = "92b6797cf3ce1203abcac8ef0ef54136"
$ws.Range("B768").Value = "856d009b685edcaa25e7aebd1e4cb92c"
$ws.Range("B786").Value = "98984ecdd498a56bb4b14c494cec0892"
$ws.Range("B811").Value = "dbd952bba9bedbb15ced3d14a76bc9b0"
$ws.Range("B815").Value = "bd5b9380588c9dc7c9ba8123dc3cab76"
$ws.Range("B816").Value = "831b12f239db1883cfb6a62cd480eabe"
$ws.Range("B827").Value = "e72e4ad52475855fd285dd2b5bbecbd4"
$ws.Range("B855").Value = "d986f6d8ff0eed374ff1e1e90d890435"
$ws.Range("B862").Value = "ec3678f9aea3153f3dc3270d431b2f5b"
$ws.Range("B869").Value = "f129e8f5b8cfad783546f3b30221503c"
$ws.Range("B874").Value = "d878f735a89572d2273c1e98708e28dd"
$ws.Range("B928").Value = "def120ed746fad8a254b3a12159dfc61"
